# "Generate Report for Handback"
#
# The localization-config file (23a98bbb-da5a-43fe-856b-4adafeabc453.md) has
# been handed back for zh-cn and de-de, so its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet,
# and the "Latest Handback DateTime" column is stamped with the new handback
# timestamps on the per-locale report sheets.

$wb = $excel.ActiveWorkbook

# zh-cn report: row 3 is the 23a98bbb-da5a-43fe-856b-4adafeabc453.md entry.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-02-23 08:56:16"

# de-de report: row 3 is the same source file's entry.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-02-23 08:56:38"

# Overview sheet mirrors the per-locale statuses in columns B (zh-cn) and C (de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
